# Trade #92 closed at 2026-02-17 21:19:02 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1401.01   # Current Capital
$summary.Range("B4").Value = 0.8       # Total P&L $
$summary.Range("B5").Value = 0.13      # Total P&L %
$summary.Range("B6").Value = 120       # Total Trades
$summary.Range("B8").Value = 47        # Losing Trades
$summary.Range("B9").Value = 44.17     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 101.01               # Capital
$status.Range("D5").Value = 87                   # Trades
$status.Range("E5").Value = 0.6899999999999999   # P&L $
$status.Range("F5").Value = 1.01                 # P&L %
$status.Range("G5").Value = 44.83                # Win Rate %

# ---------------------------------------------------------------------------
# All Trades sheet - trade #120 (row 121) closes out; trade #153 appended
# (row 154)
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Cells.Item(121, 7).Value = 0.056604        # G121 Exit Price
$allTrades.Cells.Item(121, 8).Value = "CLOSED"        # H121 Status
$allTrades.Cells.Item(121, 9).Value = -43.3962        # I121 P&L %
$allTrades.Cells.Item(121, 10).Value = -0.04          # J121 P&L $
$allTrades.Cells.Item(121, 11).Value = 101.01         # K121 Capital After
$allTrades.Cells.Item(121, 12).Value = "early_exit"   # L121 Exit Reason
$allTrades.Cells.Item(121, 13).Value = 0.15           # M121 Duration (min)

$allTrades.Cells.Item(154, 1).Value = 153
$allTrades.Cells.Item(154, 2).Value = "'2026-02-17"
$allTrades.Cells.Item(154, 3).Value = "21:18:55"
$allTrades.Cells.Item(154, 4).Value = "MarketMaking"
$allTrades.Cells.Item(154, 5).Value = "UP"
$allTrades.Cells.Item(154, 6).Value = 0.1
$allTrades.Cells.Item(154, 8).Value = "OPEN"
$allTrades.Cells.Item(154, 9).Value = 0
$allTrades.Cells.Item(154, 10).Value = 0
$allTrades.Cells.Item(154, 11).Value = 101.0503221760222
$allTrades.Cells.Item(154, 13).Value = 0
$allTrades.Cells.Item(154, 14).Value = 0
$allTrades.Cells.Item(154, 15).Value = 0
$allTrades.Cells.Item(154, 16).Value = 0.6
$allTrades.Cells.Item(154, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# MarketMaking sheet - trade #120 (row 88) closes out; trade #153 appended
# (row 121)
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

$mm.Cells.Item(88, 7).Value = 0.056604        # G88 Exit Price
$mm.Cells.Item(88, 8).Value = "CLOSED"        # H88 Status
$mm.Cells.Item(88, 9).Value = -43.3962        # I88 P&L %
$mm.Cells.Item(88, 10).Value = -0.04          # J88 P&L $
$mm.Cells.Item(88, 11).Value = 101.01         # K88 Capital After
$mm.Cells.Item(88, 16).Value = "early_exit"   # P88 Exit Reason
$mm.Cells.Item(88, 17).Value = 0.15           # Q88 Duration (min)

$mm.Cells.Item(121, 1).Value = 153
$mm.Cells.Item(121, 2).Value = "'2026-02-17"
$mm.Cells.Item(121, 3).Value = "21:18:55"
$mm.Cells.Item(121, 4).Value = "MarketMaking"
$mm.Cells.Item(121, 5).Value = "UP"
$mm.Cells.Item(121, 6).Value = 0.1
$mm.Cells.Item(121, 8).Value = "OPEN"
$mm.Cells.Item(121, 9).Value = 0
$mm.Cells.Item(121, 10).Value = 0
$mm.Cells.Item(121, 11).Value = 101.0503221760222
$mm.Cells.Item(121, 12).Value = 0
$mm.Cells.Item(121, 13).Value = 0
$mm.Cells.Item(121, 14).Value = 0.6
$mm.Cells.Item(121, 15).Value = "Normal spread capture: 19600 bps"
$mm.Cells.Item(121, 17).Value = 0
